$p = $ppt.ActivePresentation

# --- Slide 9: "Architecture – L3 {container} components" -> "Architecture – L3 statistics engine components"
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(1)

$sh9.Left = 41.511417322834646
$sh9.Top = 24.87236310472441
$sh9.Width = 876.9771653543307
$sh9.Height = 104.34330708661417

$tr9 = $sh9.TextFrame.TextRange
$tr9.Text = "Architecture – L3 statistics engine components"
# Force the same run-split boundaries as the authored edit: "Architecture – L3 " | "s" | "tatistics engine components"
$mid9 = $tr9.Characters(19, 1)
$mid9.Text = "s"

# --- Slide 10: "Architecture – L3 {container} components" -> "Architecture – L3 enrollment manager components"
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(1)

$sh10.Left = 28.645118110236222
$sh10.Top = 21.001418122834643
$sh10.Width = 922.8387401574803
$sh10.Height = 104.34330708661417

$tr10 = $sh10.TextFrame.TextRange
$tr10.Text = "Architecture – L3 enrollment manager components"
# Force the same run-split boundaries as the authored edit: "Architecture – L3 enrollment " | "manager " | "components"
$mid10 = $tr10.Characters(30, 8)
$mid10.Text = "manager "
